# "Add files via upload" — weekly progress update for 002王春妍's sheet:
# two new rows of progress notes were appended, and the active sheet/
# selection in the workbook moved from 003张保江 to 002王春妍.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("002王春妍")

# New row 16
$ws.Range("A16").Value = "20240506-20240513"
$ws.Range("B16").Value = "3.①核对ICA计算流程和数据②完成课程论文和复习备考"
$ws.Range("D16").Value = "3.准备下周的几门考试"

# New row 17 (B17 keeps the wrap-text style already used elsewhere in the sheet)
$ws.Range("A17").Value = "20240513-20240519"
$ws.Range("B17").Value = "1.编写代码`n3.①完成论文和考试②计算AD数据集的皮层厚度"
$ws.Range("B17").WrapText = $true
$ws.Range("D17").Value = "3.①准备下周的考试和课程论文"

# Make 002王春妍 the active sheet / selected tab, with the cursor on F16
# (this also clears tabSelected on whichever sheet had it before, i.e. 003张保江)
$ws.Activate()
$ws.Range("F16").Select()
